$wb = $excel.ActiveWorkbook

# --- Sheet "Feuil1": update several value cells in column C ---
$ws1 = $wb.Worksheets.Item("Feuil1")

# Plain text replacements (not number-like, safe to assign directly)
$ws1.Range("C6").Value  = "goldman sachs International"
$ws1.Range("C7").Value  = "coupon autocall"
$ws1.Range("C8").Value  = "trimestre"
$ws1.Range("C10").Value = "BNP Paribas et Stellantis NV et Veolia Environnement SA"
$ws1.Range("C11").Value = "wo action"

# Number-like / date-like text: force the cell to Text format first so the
# engine keeps it as a literal string ("1", "1.00", "100", "2023-08-02")
# instead of re-interpreting it as a number or a date serial.
$ws1.Range("C5").NumberFormat  = "@"
$ws1.Range("C9").NumberFormat  = "@"
$ws1.Range("C17").NumberFormat = "@"
$ws1.Range("C27").NumberFormat = "@"

$ws1.Range("C5").Value  = "1"
$ws1.Range("C9").Value  = "1.00"
$ws1.Range("C17").Value = "100"
$ws1.Range("C27").Value = "2023-08-02"

# --- Sheet "TRA": drop the five placeholder rows under the header row ---
$ws2 = $wb.Worksheets.Item("TRA")
$ws2.Range("A2:A6").EntireRow.Delete()

# --- Sheet "DATE": refresh the quarterly payment/redemption schedule ---
$ws3 = $wb.Worksheets.Item("DATE")
$ws3.Range("A2").Value = "31/07/2023, 30/10/2023, 29/01/2024, 29/04/2024, 29/07/2024, 29/10/2024, 29/01/2025, 29/04/2025, 29/07/2025, 29/10/2025, 29/01/2026, 29/04/2026, 29/07/2026, 29/10/2026, 29/01/2027, 29/04/2027, 29/07/2027, 29/10/2027, 31/01/2028, 02/05/2028, 31/07/2028, 30/10/2028, 29/01/2029, 30/04/2029, 30/07/2029, 29/10/2029, 29/01/2030, 29/04/2030, 29/07/2030, 29/10/2030, 29/01/2031, 29/04/2031, 29/07/2031, 29/10/2031, 29/01/2032, 29/04/2032, 29/07/2032, 29/07/2032"
$ws3.Range("A3").Value = "Dates de remboursement"
$ws3.Range("A4").Value = "07/08/2023, 06/11/2023, 05/02/2024, 07/05/2024, 05/08/2024, 05/11/2024, 05/02/2025, 07/05/2025, 05/08/2025, 05/11/2025, 05/02/2026, 07/05/2026, 05/08/2026, 05/11/2026, 05/02/2027, 06/05/2027, 05/08/2027, 05/11/2027, 07/02/2028, 09/05/2028, 07/08/2028, 06/11/2028, 05/02/2029, 08/05/2029, 06/08/2029, 05/11/2029, 05/02/2030, 07/05/2030, 05/08/2030, 05/11/2030, 05/02/2031, 07/05/2031, 05/08/2031, 05/11/2031, 05/02/2032, 06/05/2032"
$ws3.Range("A5:A6").EntireRow.Delete()
